# ---------------------------------------------------------------------------
# Target diff analysis
# ---------------------------------------------------------------------------
# The supplied unified diff touches four parts: word/document.xml,
# word/footer.xml, word/header.xml and word/styles.xml. In every one of
# those four files the *only* change is that the root element's xmlns:*
# attributes are re-ordered (e.g. xmlns:r first vs. xmlns:w first). That is
# a side effect of which XML serializer/JAXB run produced the package - it
# is not a document edit, carries no semantic meaning, and is not
# reachable (or meaningful) through the Word object model.
#
# The single other change, in word/document.xml, is inside an XML comment
# that docx4j places as the very first child of <w:body>:
#
#   <!-- Modified by docx4j 11.5.6 (Apache licensed) using REFERENCE JAXB
#        in Microsoft Java 21.0.8 on Mac OS X -->
#            -> "Oracle Java 21.0.8"
#
# This is a build-tool stamp recording which JDK vendor ran docx4j when the
# test fixture was regenerated; it is injected by docx4j itself, not typed
# by a user in Word. It is a raw XML comment sitting directly under
# <w:body> - not body text, not a field, not a reviewing/Comments-pane
# comment (Document.Comments / Range.Comments), and not Custom XML
# mark-up (Document.XMLNodes / selectNodes / selectSingleNode), so none of
# the Word automation surfaces (Find & Replace, Range.Text, Content.XML /
# WordOpenXML, XMLNode navigation, etc.) can see or touch it - confirmed
# empirically against this host (Content.Find.Execute does not match it,
# Content.WordOpenXML/Document.XML do not include it, selectSingleNode
# ("//comment()") returns null). On real Word/COM the same is true: that
# text simply isn't part of the document model Word exposes, which is
# consistent with it only ever being written by the docx4j library itself,
# not by editing the file in Word.
#
# So there is nothing in this diff that a Word COM-interop script can
# legitimately perform. To stay faithful to the document, this script
# makes no structural edits (touching the body would, per testing, even
# risk *dropping* that already-harmless comment when Word/the host
# re-serializes the paragraph it sits next to). It only makes one
# best-effort, side-effect-free attempt at the textual substitution in
# case the wording is ever exposed as plain content; Find.Execute reports
# "not found" here and leaves the document untouched, which is the
# correct, closest achievable result for this particular diff.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$d.Content.Find.Execute("Microsoft Java 21.0.8", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Oracle Java 21.0.8", 2) | Out-Null
